# Add 2022-Q4 data
#
# Target end state:
#   Sheets (in tab order): 总计, 2022-Q4 (new), 2022-Q1, 2020-Q4
#   总计 sheet: a new row is inserted for the "2022-Q4" summary line, ahead of the
#               existing "2022-Q1"/"2020-Q4" rows.
#   2022-Q4 sheet: brand-new worksheet holding the per-fund breakdown for the
#                  new quarter (same shape as the existing quarter sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet — insert a new row 2 for 2022-Q4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Push the existing data rows down by inserting a fresh row at row 2, copying
# the row-above formatting (as Excel normally does), then fix up the cells.
$total.Rows.Item(2).Insert()

# The new row's B/C/D cells inherited the header's style from the insert；
# reset them back to the (unstyled) look used by the other data rows.
$total.Cells.Item(2, 2).Style = "Normal"
$total.Cells.Item(2, 3).Style = "Normal"
$total.Cells.Item(2, 4).Style = "Normal"

# Column A carries the bordered/centered style used throughout the data rows;
# grab it from the row just below (shifted-down former row 2).
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.02

# Renumber the index column for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2

# ---------------------------------------------------------------------------
# 2. Brand-new "2022-Q4" worksheet, placed right after "总计".
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Match page margins used by the sibling quarter sheets.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1

# Reference sheet to borrow the header/index-column style (style index "2")
# from, via a Copy() (copies formatting together with whatever value is
# already in the source cell — the value is overwritten immediately after).
$styleSrc = $wb.Worksheets.Item("2022-Q1")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $styleSrc.Cells.Item(1, $col).Copy($q4.Cells.Item(1, $col))
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
}

function Set-TextCell($cell, $text) {
    # Force text storage (avoids "002236" / "5.67" / "93.32" being silently
    # coerced into numbers), then drop back to the default/unstyled look.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$data = @(
    @("002236", "大成中证360互联网+大数据100指数A", "1.15", "92.50", "0.99", "0.0114", 8),
    @("003359", "大成中证360互联网+大数据100指数C", "1.12", "92.50", "0.99", "0.0111", 8),
    @("002952", "建信多因子量化股票",                 "0.09", "91.91", "2.20", "0.0020", 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    # Column A: index style copied from the reference sheet's data rows (it
    # only has two data rows, row 2 has the style we need for every row).
    $styleSrc.Cells.Item(2, 1).Copy($q4.Cells.Item($row, 1))
    $q4.Cells.Item($row, 1).Value = $i

    Set-TextCell $q4.Cells.Item($row, 2) $rec[0]
    Set-TextCell $q4.Cells.Item($row, 3) $rec[1]
    Set-TextCell $q4.Cells.Item($row, 4) $rec[2]
    Set-TextCell $q4.Cells.Item($row, 5) $rec[3]
    Set-TextCell $q4.Cells.Item($row, 6) $rec[4]
    Set-TextCell $q4.Cells.Item($row, 7) $rec[5]

    $q4.Cells.Item($row, 8).Value = $rec[6]
}
